# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 53
$ws1.Range("F7").Value = 2628
$ws1.Range("F8").Value = 1145
$ws1.Range("F9").Value = 230
$ws1.Range("F10").Value = 92
$ws1.Range("F11").Value = 5902
$ws1.Range("F12").Value = 69
$ws1.Range("F14").Value = 578
$ws1.Range("F15").Value = 11574
$ws1.Range("F16").Value = 11759

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 53
$ws4.Range("F7").Value = 2628
$ws4.Range("F9").Value = 1145
$ws4.Range("F10").Value = 230
$ws4.Range("F11").Value = 92
$ws4.Range("F12").Value = 5902
$ws4.Range("F13").Value = 69
$ws4.Range("F15").Value = 578
$ws4.Range("F16").Value = 11574
$ws4.Range("F17").Value = 11759
